$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 224.7833046666667
$ws.Range("H2").Value = 674.349914
$ws.Range("I2").Value = 0.3882379172278888
$ws.Range("J2").Value = 0.3882379172278889
$ws.Range("M2").Value = 3.556762333333333
$ws.Range("N2").Value = 10.670287
$ws.Range("O2").Value = 0.04280930450251701
$ws.Range("P2").Value = 0.04280930450251701
$ws.Range("Q2").Value = 799.5007912005909
$ws.Range("R2").Value = 7195.507120805318
$ws.Range("S2").Value = 0.01662019521803169
$ws.Range("T2").Value = 0.01662019521803169
$ws.Range("G3").Value = 224.7833046666667
$ws.Range("H3").Value = 674.349914
$ws.Range("I3").Value = 0.3882379172278888
$ws.Range("J3").Value = 0.3882379172278889
$ws.Range("O3").Value = 0.5686906263805706
$ws.Range("P3").Value = 0.5686906263805704
$ws.Range("Q3").Value = 10620.78935930606
$ws.Range("R3").Value = 95587.10423375457
$ws.Range("S3").Value = 0.2207872643330162
$ws.Range("T3").Value = 0.2207872643330162
$ws.Range("G4").Value = 224.7833046666667
$ws.Range("H4").Value = 674.349914
$ws.Range("I4").Value = 0.3882379172278888
$ws.Range("J4").Value = 0.3882379172278889
$ws.Range("M4").Value = 24.53173066666666
$ws.Range("N4").Value = 73.595192
$ws.Range("O4").Value = 0.2952646900921413
$ws.Range("P4").Value = 0.2952646900921412
$ws.Range("Q4").Value = 5514.323488445943
$ws.Range("R4").Value = 49628.91139601349
$ws.Range("S4").Value = 0.114632948312311
$ws.Range("T4").Value = 0.114632948312311
$ws.Range("G5").Value = 224.7833046666667
$ws.Range("H5").Value = 674.349914
$ws.Range("I5").Value = 0.3882379172278888
$ws.Range("J5").Value = 0.3882379172278889
$ws.Range("M5").Value = 7.746355333333334
$ws.Range("N5").Value = 23.239066
$ws.Range("O5").Value = 0.09323537902477132
$ws.Range("P5").Value = 0.0932353790247713
$ws.Range("Q5").Value = 1741.251350948925
$ws.Range("R5").Value = 15671.26215854032
$ws.Range("S5").Value = 0.03619750936453001
$ws.Range("T5").Value = 0.03619750936453001
$ws.Range("I6").Value = 0.4251955538547045
$ws.Range("J6").Value = 0.4251955538547046
$ws.Range("M6").Value = 3.556762333333333
$ws.Range("N6").Value = 10.670287
$ws.Range("O6").Value = 0.04280930450251701
$ws.Range("P6").Value = 0.04280930450251701
$ws.Range("Q6").Value = 875.6078853634187
$ws.Range("R6").Value = 7880.470968270769
$ws.Range("S6").Value = 0.01820232593808242
$ws.Range("T6").Value = 0.01820232593808242
$ws.Range("I7").Value = 0.4251955538547045
$ws.Range("J7").Value = 0.4251955538547046
$ws.Range("O7").Value = 0.5686906263805706
$ws.Range("P7").Value = 0.5686906263805704
$ws.Range("S7").Value = 0.2418047258558655
$ws.Range("T7").Value = 0.2418047258558655
$ws.Range("I8").Value = 0.4251955538547045
$ws.Range("J8").Value = 0.4251955538547046
$ws.Range("M8").Value = 24.53173066666666
$ws.Range("N8").Value = 73.595192
$ws.Range("O8").Value = 0.2952646900921413
$ws.Range("P8").Value = 0.2952646900921412
$ws.Range("Q8").Value = 6039.249969568277
$ws.Range("R8").Value = 54353.2497261145
$ws.Range("S8").Value = 0.1255452334374657
$ws.Range("T8").Value = 0.1255452334374657
$ws.Range("I9").Value = 0.4251955538547045
$ws.Range("J9").Value = 0.4251955538547046
$ws.Range("M9").Value = 7.746355333333334
$ws.Range("N9").Value = 23.239066
$ws.Range("O9").Value = 0.09323537902477132
$ws.Range("P9").Value = 0.0932353790247713
$ws.Range("Q9").Value = 1907.006759807016
$ws.Range("R9").Value = 17163.06083826314
$ws.Range("S9").Value = 0.03964326862329094
$ws.Range("T9").Value = 0.03964326862329094
$ws.Range("G10").Value = 107.695137
$ws.Range("H10").Value = 323.085411
$ws.Range("I10").Value = 0.186007300437435
$ws.Range("J10").Value = 0.186007300437435
$ws.Range("M10").Value = 3.556762333333333
$ws.Range("N10").Value = 10.670287
$ws.Range("O10").Value = 0.04280930450251701
$ws.Range("P10").Value = 0.04280930450251701
$ws.Range("Q10").Value = 383.046006764773
$ws.Range("R10").Value = 3447.414060882957
$ws.Range("S10").Value = 0.007962843164117319
$ws.Range("T10").Value = 0.007962843164117321
$ws.Range("G11").Value = 107.695137
$ws.Range("H11").Value = 323.085411
$ws.Range("I11").Value = 0.186007300437435
$ws.Range("J11").Value = 0.186007300437435
$ws.Range("O11").Value = 0.5686906263805706
$ws.Range("P11").Value = 0.5686906263805704
$ws.Range("Q11").Value = 5088.48896404816
$ws.Range("R11").Value = 45796.40067643344
$ws.Range("S11").Value = 0.1057806081971239
$ws.Range("T11").Value = 0.1057806081971239
$ws.Range("G12").Value = 107.695137
$ws.Range("H12").Value = 323.085411
$ws.Range("I12").Value = 0.186007300437435
$ws.Range("J12").Value = 0.186007300437435
$ws.Range("M12").Value = 24.53173066666666
$ws.Range("N12").Value = 73.595192
$ws.Range("O12").Value = 0.2952646900921413
$ws.Range("P12").Value = 0.2952646900921412
$ws.Range("Q12").Value = 2641.948094993768
$ws.Range("R12").Value = 23777.53285494391
$ws.Range("S12").Value = 0.05492138791853505
$ws.Range("T12").Value = 0.05492138791853506
$ws.Range("G13").Value = 107.695137
$ws.Range("H13").Value = 323.085411
$ws.Range("I13").Value = 0.186007300437435
$ws.Range("J13").Value = 0.186007300437435
$ws.Range("M13").Value = 7.746355333333334
$ws.Range("N13").Value = 23.239066
$ws.Range("O13").Value = 0.09323537902477132
$ws.Range("P13").Value = 0.0932353790247713
$ws.Range("Q13").Value = 834.244798874014
$ws.Range("R13").Value = 7508.203189866127
$ws.Range("S13").Value = 0.01734246115765876
$ws.Range("T13").Value = 0.01734246115765876
$ws.Range("G14").Value = 0.323784
$ws.Range("H14").Value = 0.971352
$ws.Range("I14").Value = 0.0005592284799715185
$ws.Range("J14").Value = 0.0005592284799715186
$ws.Range("M14").Value = 3.556762333333333
$ws.Range("N14").Value = 10.670287
$ws.Range("O14").Value = 0.04280930450251701
$ws.Range("P14").Value = 0.04280930450251701
$ws.Range("Q14").Value = 1.151622735336
$ws.Range("R14").Value = 10.364604618024
$ws.Range("S14").Value = 0.00002394018228558047
$ws.Range("T14").Value = 0.00002394018228558047
$ws.Range("G15").Value = 0.323784
$ws.Range("H15").Value = 0.971352
$ws.Range("I15").Value = 0.0005592284799715185
$ws.Range("J15").Value = 0.0005592284799715186
$ws.Range("O15").Value = 0.5686906263805706
$ws.Range("P15").Value = 0.5686906263805704
$ws.Range("Q15").Value = 15.29847453312
$ws.Range("R15").Value = 137.68627079808
$ws.Range("S15").Value = 0.0003180279945648572
$ws.Range("T15").Value = 0.0003180279945648572
$ws.Range("G16").Value = 0.323784
$ws.Range("H16").Value = 0.971352
$ws.Range("I16").Value = 0.0005592284799715185
$ws.Range("J16").Value = 0.0005592284799715186
$ws.Range("M16").Value = 24.53173066666666
$ws.Range("N16").Value = 73.595192
$ws.Range("O16").Value = 0.2952646900921413
$ws.Range("P16").Value = 0.2952646900921412
$ws.Range("Q16").Value = 7.942981882175999
$ws.Range("R16").Value = 71.48683693958399
$ws.Range("S16").Value = 0.0001651204238294897
$ws.Range("T16").Value = 0.0001651204238294897
$ws.Range("G17").Value = 0.323784
$ws.Range("H17").Value = 0.971352
$ws.Range("I17").Value = 0.0005592284799715185
$ws.Range("J17").Value = 0.0005592284799715186
$ws.Range("M17").Value = 7.746355333333334
$ws.Range("N17").Value = 23.239066
$ws.Range("O17").Value = 0.09323537902477132
$ws.Range("P17").Value = 0.0932353790247713
$ws.Range("Q17").Value = 2.508145915248
$ws.Range("R17").Value = 22.573313237232
$ws.Range("S17").Value = 0.00005213987929159126
$ws.Range("T17").Value = 0.00005213987929159127
